# Update workbook with Tanzania data.
# - Switch the "filter programs" objective from "thrive" to "thrive, child_deaths"
# - Update the set of programs marked "included" on the "Programs to include" sheet
# - Update the matching coverage targets on the "Coverage scenario" sheet
# - Move the active selection/tab from "Optimisation options" to "Programs to include"

$wb = $excel.ActiveWorkbook

$wsPrograms  = $wb.Worksheets.Item("Programs to include")
$wsCoverage  = $wb.Worksheets.Item("Coverage scenario")
$wsOptim     = $wb.Worksheets.Item("Optimisation options")

# --- Programs to include: flip which programs are marked with "x" ---------

# Newly included programs
$wsPrograms.Range("B5").Value  = "x"   # Cash transfers
$wsPrograms.Range("B7").Value  = "x"   # IFA fortification of maize
$wsPrograms.Range("B14").Value = "x"   # IFAS for pregnant women (community)
$wsPrograms.Range("B16").Value = "x"   # IPTp

# No-longer-included programs
$wsPrograms.Range("B31").Value = ""    # WASH: Handwashing
$wsPrograms.Range("B32").Value = ""    # WASH: Hygenic disposal
$wsPrograms.Range("B33").Value = ""    # WASH: Improved sanitation
$wsPrograms.Range("B34").Value = ""    # WASH: Improved water source
$wsPrograms.Range("B35").Value = ""    # WASH: Piped water
$wsPrograms.Range("B37").Value = ""    # Zinc supplementation

# --- Coverage scenario: update 2018 (column D) coverage targets -----------

$wsCoverage.Range("D14").Value = 0.95  # IFAS for pregnant women (community)

$wsCoverage.Range("D25").Value = ""    # Micronutrient powders
$wsCoverage.Range("D31").Value = ""    # WASH: Handwashing
$wsCoverage.Range("D32").Value = ""    # WASH: Hygenic disposal
$wsCoverage.Range("D33").Value = ""    # WASH: Improved sanitation
$wsCoverage.Range("D34").Value = ""    # WASH: Improved water source
$wsCoverage.Range("D35").Value = ""    # WASH: Piped water
$wsCoverage.Range("D36").Value = ""    # Zinc for treatment + ORS
$wsCoverage.Range("D37").Value = ""    # Zinc supplementation

# --- Optimisation options: update objectives -------------------------------

$wsOptim.Range("C2").Value = "thrive, child_deaths"

# --- Selection / active sheet bookkeeping ----------------------------------

$wsCoverage.Range("D14").Select() | Out-Null
$wsOptim.Range("C2").Select() | Out-Null
$wsPrograms.Range("B15").Select() | Out-Null
$wsPrograms.Select() | Out-Null
